$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp (A1)
$ws.Range('A1').Value = 'Datos actualizados a 20 de Abril de 2020 a las 00:52'

# Row 4: Estados Unidos
$ws.Range('A4').Value = 'Estados Unidos'
$ws.Range('B4').Value = 763083
$ws.Range('C4').Value = 24291
$ws.Range('D4').Value = 70806
$ws.Range('E4').Value = 651782
$ws.Range('F4').Value = 13566
$ws.Range('G4').Value = 1481
$ws.Range('H4').Value = 40495

# Row 5: España
$ws.Range('A5').Value = 'España'
$ws.Range('B5').Value = 198674
$ws.Range('C5').Value = 4258
$ws.Range('D5').Value = 77357
$ws.Range('E5').Value = 100864
$ws.Range('F5').Value = 7371
$ws.Range('G5').Value = 410
$ws.Range('H5').Value = 20453

# Row 14: Brasil
$ws.Range('A14').Value = 'Brasil'
$ws.Range('B14').Value = 38654
$ws.Range('C14').Value = 1932
$ws.Range('D14').Value = 22130
$ws.Range('E14').Value = 14062
$ws.Range('F14').Value = 6634
$ws.Range('G14').Value = 101
$ws.Range('H14').Value = 2462

# Row 16: Canada
$ws.Range('A16').Value = 'Canada'
$ws.Range('B16').Value = 35056
$ws.Range('C16').Value = 1673
$ws.Range('D16').Value = 11843
$ws.Range('E16').Value = 21626
$ws.Range('F16').Value = 557
$ws.Range('G16').Value = 117
$ws.Range('H16').Value = 1587

# Row 38: Chequia
$ws.Range('A38').Value = 'Chequia'
$ws.Range('B38').Value = 6746
$ws.Range('C38').Value = 140
$ws.Range('D38').Value = 1298
$ws.Range('E38').Value = 5262
$ws.Range('F38').Value = 84
$ws.Range('G38').Value = 5
$ws.Range('H38').Value = 186

# Row 50: Colombia
$ws.Range('A50').Value = 'Colombia'
$ws.Range('B50').Value = 3792
$ws.Range('C50').Value = 171
$ws.Range('D50').Value = 711
$ws.Range('E50').Value = 2902
$ws.Range('F50').Value = 98
$ws.Range('G50').Value = 13
$ws.Range('H50').Value = 179

# Row 51: Finlandia
$ws.Range('A51').Value = 'Finlandia'
$ws.Range('B51').Value = 3783
$ws.Range('C51').Value = 102
$ws.Range('D51').Value = 1700
$ws.Range('E51').Value = 1989
$ws.Range('F51').Value = 68
$ws.Range('G51').Value = 4
$ws.Range('H51').Value = 94

# Row 80: Ghana
$ws.Range('A80').Value = 'Ghana'
$ws.Range('B80').Value = 1042
$ws.Range('C80').Value = 208
$ws.Range('D80').Value = 99
$ws.Range('E80').Value = 934
$ws.Range('F80').Value = 4
$ws.Range('G80').Value = 0
$ws.Range('H80').Value = 9

# Row 81: Cuba
$ws.Range('A81').Value = 'Cuba'
$ws.Range('B81').Value = 1035
$ws.Range('C81').Value = 49
$ws.Range('D81').Value = 255
$ws.Range('E81').Value = 746
$ws.Range('F81').Value = 16
$ws.Range('G81').Value = 2
$ws.Range('H81').Value = 34

# Row 82: Hong Kong
$ws.Range('A82').Value = 'Hong Kong'
$ws.Range('B82').Value = 1026
$ws.Range('C82').Value = 2
$ws.Range('D82').Value = 602
$ws.Range('E82').Value = 420
$ws.Range('F82').Value = 8
$ws.Range('G82').Value = 0
$ws.Range('H82').Value = 4

# Row 83: Camerun
$ws.Range('A83').Value = 'Camerun'
$ws.Range('B83').Value = 1017
$ws.Range('C83').Value = 0
$ws.Range('D83').Value = 305
$ws.Range('E83').Value = 670
$ws.Range('F83').Value = 33
$ws.Range('G83').Value = 0
$ws.Range('H83').Value = 42

# Row 84: Afganistan
$ws.Range('A84').Value = 'Afganistan'
$ws.Range('B84').Value = 996
$ws.Range('C84').Value = 63
$ws.Range('D84').Value = 131
$ws.Range('E84').Value = 832
$ws.Range('F84').Value = 7
$ws.Range('G84').Value = 3
$ws.Range('H84').Value = 33

# Row 85: Bulgaria
$ws.Range('A85').Value = 'Bulgaria'
$ws.Range('B85').Value = 894
$ws.Range('C85').Value = 16
$ws.Range('D85').Value = 161
$ws.Range('E85').Value = 691
$ws.Range('F85').Value = 36
$ws.Range('G85').Value = 1
$ws.Range('H85').Value = 42

# Row 86: Tunez
$ws.Range('A86').Value = 'Tunez'
$ws.Range('B86').Value = 879
$ws.Range('C86').Value = 13
$ws.Range('D86').Value = 43
$ws.Range('E86').Value = 798
$ws.Range('F86').Value = 33
$ws.Range('G86').Value = 1
$ws.Range('H86').Value = 38

# Row 87: Costa de Marfil
$ws.Range('A87').Value = 'Costa de Marfil'
$ws.Range('B87').Value = 847
$ws.Range('C87').Value = 46
$ws.Range('D87').Value = 260
$ws.Range('E87').Value = 578
$ws.Range('F87').Value = 0
$ws.Range('G87').Value = 1
$ws.Range('H87').Value = 9

# Row 88: Republica de Yibuti
$ws.Range('A88').Value = 'Republica de Yibuti'
$ws.Range('B88').Value = 846
$ws.Range('C88').Value = 114
$ws.Range('D88').Value = 102
$ws.Range('E88').Value = 742
$ws.Range('F88').Value = 0
$ws.Range('G88').Value = 0
$ws.Range('H88').Value = 2

# Row 94: Costa Rica
$ws.Range('A94').Value = 'Costa Rica'
$ws.Range('B94').Value = 660
$ws.Range('C94').Value = 5
$ws.Range('D94').Value = 112
$ws.Range('E94').Value = 543
$ws.Range('F94').Value = 10
$ws.Range('G94').Value = 1
$ws.Range('H94').Value = 5
